$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (Beta) values ---
$ws.Range("C2").Value = 40.35273080373261
$ws.Range("E2").Value = 0.05768901818751787
$ws.Range("F2").Value = 40.32575735259291
$ws.Range("G2").Value = 39.69290616083104
$ws.Range("H2").Value = 40.93257907186761
$ws.Range("I2").Value = 0.0007613987687944709
$ws.Range("J2").Value = 0.0007110440303010838
$ws.Range("K2").Value = 0.0008537742827225605
$ws.Range("L2").Value = 0.05766958104364201
$ws.Range("M2").Value = 0.0572362877039788
$ws.Range("N2").Value = 0.05810820245610884

# --- Update existing row 3 (Gamma) values ---
$ws.Range("F3").Value = 0.00001480883976174689
$ws.Range("G3").Value = 0.000000007401131980785392
$ws.Range("H3").Value = 0.00004180993561973118
$ws.Range("I3").Value = 0.00001298066649109086
$ws.Range("J3").Value = 0.000000006918286695062993
$ws.Range("K3").Value = 0.00003648888266881516
$ws.Range("L3").Value = 0.00001525364180330135
$ws.Range("M3").Value = 0.000000007646185072349991
$ws.Range("N3").Value = 0.00004304144698298873

# --- Add new row 4 (Beta + Gamma) ---
$ws.Range("A4").Value = 2
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 40.35273080373261
$ws.Range("D4").Value = 0.0007432820064133916
$ws.Range("E4").Value = 0.05768901818751787
$ws.Range("F4").Value = 40.32577216143267
$ws.Range("G4").Value = 39.69290616823217
$ws.Range("H4").Value = 40.93262088180324
$ws.Range("I4").Value = 0.0007743794352855617
$ws.Range("J4").Value = 0.0007110509485877788
$ws.Range("K4").Value = 0.0008902631653913758
$ws.Range("L4").Value = 0.05768483468544532
$ws.Range("M4").Value = 0.05723629535016387
$ws.Range("N4").Value = 0.05815124390309181
